# Auto-generated Excel COM-interop script
# - Appends 44 event rows to sheet '展览' (Exhibition)
# - Appends 1 event row to sheet '演出' (Performance)
# - Fills in the '想去人数' (F) column on '全部类型' (All types),
#   which previously held placeholder zeros.

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 (Exhibition) -- append data rows 2..45 ----
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("B2:B45").NumberFormat = "@"

$arrExpo = New-Object 'object[,]' 44,9
$arrExpo[0,0] = 1
$arrExpo[0,1] = '2024-06-30'
$arrExpo[0,2] = '南昌·ChinastyleCOSPLAY  '
$arrExpo[0,3] = '真君路888号 南昌华侨城玩美公园'
$arrExpo[0,4] = '2024.06.30 09:30-07.02 17:30'
$arrExpo[0,5] = 1862
$arrExpo[0,6] = 75
$arrExpo[0,7] = 'https://show.bilibili.com/platform/detail.html?id=87045'
$arrExpo[0,8] = '//i1.hdslb.com/bfs/openplatform/202406/OEU3ijdb1719299094349.jpeg'
$arrExpo[1,0] = 2
$arrExpo[1,1] = '2024-06-30'
$arrExpo[1,2] = '宜春·BM次元盛典运动番only（取消）'
$arrExpo[1,3] = '鼓楼西路与官圳路交叉口东120米 地中海宴会酒店(润达店)'
$arrExpo[1,4] = '2024.06.30 10:00-06.30 17:00'
$arrExpo[1,5] = 277
$arrExpo[1,6] = '不可售'
$arrExpo[1,7] = 'https://show.bilibili.com/platform/detail.html?id=84636'
$arrExpo[1,8] = '//i1.hdslb.com/bfs/openplatform/202405/oaGZXKok1715328213440.png'
$arrExpo[2,0] = 3
$arrExpo[2,1] = '2024-07-06'
$arrExpo[2,2] = '南昌·次元星球动漫游戏展'
$arrExpo[2,3] = '龙蟠街666号融创茂1层 融创茂'
$arrExpo[2,4] = '2024.07.06 10:00-07.06 17:00'
$arrExpo[2,5] = 32
$arrExpo[2,6] = '不可售'
$arrExpo[2,7] = 'https://show.bilibili.com/platform/detail.html?id=86405'
$arrExpo[2,8] = '//i2.hdslb.com/bfs/openplatform/202405/9ZfGuXJ01716796674559.jpeg'
$arrExpo[3,0] = 4
$arrExpo[3,1] = '2024-07-06'
$arrExpo[3,2] = '鹰潭·BM次元盛典运动番only（取消）'
$arrExpo[3,3] = '体育馆东路2号九小隔壁 忆江南•宴会楼'
$arrExpo[3,4] = '2024.07.06 10:00-07.06 17:00'
$arrExpo[3,5] = 62
$arrExpo[3,6] = '不可售'
$arrExpo[3,7] = 'https://show.bilibili.com/platform/detail.html?id=85997'
$arrExpo[3,8] = '//i1.hdslb.com/bfs/openplatform/202405/4yuR8NQc1716259522268.png'
$arrExpo[4,0] = 5
$arrExpo[4,1] = '2024-07-07'
$arrExpo[4,2] = '赣州·BM次元盛典运动番only（取消）'
$arrExpo[4,3] = '米瑞金路2口0号上客天下1楼 上客天下.老虔州'
$arrExpo[4,4] = '2024.07.07 10:00-07.07 17:00'
$arrExpo[4,5] = 44
$arrExpo[4,6] = '不可售'
$arrExpo[4,7] = 'https://show.bilibili.com/platform/detail.html?id=86602'
$arrExpo[4,8] = '//i1.hdslb.com/bfs/openplatform/202405/Xrq9sfkE1716259438090.png'
$arrExpo[5,0] = 6
$arrExpo[5,1] = '2024-07-12'
$arrExpo[5,2] = '新余·2024第三届MG动漫嘉年华'
$arrExpo[5,3] = '仙女湖大道与五一南路交叉口西约180米 老上海风情街水晶厅'
$arrExpo[5,4] = '2024.07.12 10:00-07.13 17:30'
$arrExpo[5,5] = 169
$arrExpo[5,6] = 55
$arrExpo[5,7] = 'https://show.bilibili.com/platform/detail.html?id=86536'
$arrExpo[5,8] = '//i0.hdslb.com/bfs/openplatform/202405/11RbfeFq1716813676323.jpeg'
$arrExpo[6,0] = 7
$arrExpo[6,1] = '2024-07-13'
$arrExpo[6,2] = '南昌·SuperComic动漫游戏博览会'
$arrExpo[6,3] = '怀玉山大道1315号 南昌绿地国际博览中心'
$arrExpo[6,4] = '2024.07.13 09:00-07.14 17:00'
$arrExpo[6,5] = 3700
$arrExpo[6,6] = 65
$arrExpo[6,7] = 'https://show.bilibili.com/platform/detail.html?id=86992'
$arrExpo[6,8] = '//i1.hdslb.com/bfs/openplatform/202406/wQTAjelJ1717642148929.jpeg'
$arrExpo[7,0] = 8
$arrExpo[7,1] = '2024-07-13'
$arrExpo[7,2] = '南昌·SuperComic配音演员刘明月专场见面会'
$arrExpo[7,3] = '怀玉山大道1315号 南昌绿地国际博览中心'
$arrExpo[7,4] = '2024.07.13 09:00-07.13 17:00'
$arrExpo[7,5] = 165
$arrExpo[7,6] = 168
$arrExpo[7,7] = 'https://show.bilibili.com/platform/detail.html?id=87570'
$arrExpo[7,8] = '//i2.hdslb.com/bfs/openplatform/202406/1D1reIl81718609013880.png'
$arrExpo[8,0] = 9
$arrExpo[8,1] = '2024-07-13'
$arrExpo[8,2] = '南昌·THO-梦违赣鄱荟萃·叁~幻想Strawberry~!!'
$arrExpo[8,3] = '民德路411号 东方豪景花园酒店(民德路店)'
$arrExpo[8,4] = '2024.07.13 09:30-07.13 17:30'
$arrExpo[8,5] = 108
$arrExpo[8,6] = 65
$arrExpo[8,7] = 'https://show.bilibili.com/platform/detail.html?id=87668'
$arrExpo[8,8] = '//i1.hdslb.com/bfs/openplatform/202406/Bk9cYryT1718360290362.jpeg'
$arrExpo[9,0] = 10
$arrExpo[9,1] = '2024-07-13'
$arrExpo[9,2] = '宜春·COMIC WORLD次元创作同人季·动漫游戏嘉年华'
$arrExpo[9,3] = '宜春国际商贸城会展中心 宜春国际商贸城会展中心'
$arrExpo[9,4] = '2024.07.13 10:00-07.14 17:00'
$arrExpo[9,5] = 84
$arrExpo[9,6] = 55
$arrExpo[9,7] = 'https://show.bilibili.com/platform/detail.html?id=86667'
$arrExpo[9,8] = '//i2.hdslb.com/bfs/openplatform/202405/JEjmQOLw1716737193284.jpeg'
$arrExpo[10,0] = 11
$arrExpo[10,1] = '2024-07-13'
$arrExpo[10,2] = '赣州·十万伏特-次元交流会（夏）'
$arrExpo[10,3] = '梅关大道36-16号 麋鹿星球艺术中心'
$arrExpo[10,4] = '2024.07.13 09:30-07.13 17:00'
$arrExpo[10,5] = 67
$arrExpo[10,6] = 45
$arrExpo[10,7] = 'https://show.bilibili.com/platform/detail.html?id=87597'
$arrExpo[10,8] = '//i1.hdslb.com/bfs/openplatform/202406/87yQ4Hmf1718681348727.jpeg'
$arrExpo[11,0] = 12
$arrExpo[11,1] = '2024-07-14'
$arrExpo[11,2] = '南昌·赛马娘ONLY'
$arrExpo[11,3] = '洪城路99号 锦都皇冠酒店(八一广场火车站店)'
$arrExpo[11,4] = '2024.07.14 09:00-07.14 17:30'
$arrExpo[11,5] = 72
$arrExpo[11,6] = 68
$arrExpo[11,7] = 'https://show.bilibili.com/platform/detail.html?id=87367'
$arrExpo[11,8] = '//i1.hdslb.com/bfs/openplatform/202406/wXQuIKtu1718165450704.png'
$arrExpo[12,0] = 13
$arrExpo[12,1] = '2024-07-14'
$arrExpo[12,2] = '吉安·COMIC LIFE次元假日05'
$arrExpo[12,3] = '东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心'
$arrExpo[12,4] = '2024.07.14 09:00-07.14 18:00'
$arrExpo[12,5] = 646
$arrExpo[12,6] = 52.1
$arrExpo[12,7] = 'https://show.bilibili.com/platform/detail.html?id=85924'
$arrExpo[12,8] = '//i2.hdslb.com/bfs/openplatform/202405/tBNLb2671716182857904.jpeg'
$arrExpo[13,0] = 14
$arrExpo[13,1] = '2024-07-19'
$arrExpo[13,2] = '九江·第一届Loading加载中动漫展'
$arrExpo[13,3] = '湓浦街道大中路339号 百嘉洲际酒店'
$arrExpo[13,4] = '2024.07.19 09:00-07.21 17:00'
$arrExpo[13,5] = 146
$arrExpo[13,6] = 36.6
$arrExpo[13,7] = 'https://show.bilibili.com/platform/detail.html?id=87787'
$arrExpo[13,8] = '//i0.hdslb.com/bfs/openplatform/202406/TH3lVD5G1718158901239.jpeg'
$arrExpo[14,0] = 15
$arrExpo[14,1] = '2024-07-19'
$arrExpo[14,2] = '赣州·第四届赣州半夏动漫展'
$arrExpo[14,3] = '105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心'
$arrExpo[14,4] = '2024.07.19 10:00-07.21 17:00'
$arrExpo[14,5] = 841
$arrExpo[14,6] = 55
$arrExpo[14,7] = 'https://show.bilibili.com/platform/detail.html?id=86587'
$arrExpo[14,8] = '//i1.hdslb.com/bfs/openplatform/202405/tlfL9oq91717053081587.jpeg'
$arrExpo[15,0] = 16
$arrExpo[15,1] = '2024-07-20'
$arrExpo[15,2] = '九江·第五届ACD动漫游戏嘉年华'
$arrExpo[15,3] = '九瑞大道与重庆路交汇处西南角 九江国际会展中心'
$arrExpo[15,4] = '2024.07.20 10:00-07.21 17:00'
$arrExpo[15,5] = 41
$arrExpo[15,6] = 39.9
$arrExpo[15,7] = 'https://show.bilibili.com/platform/detail.html?id=88221'
$arrExpo[15,8] = '//i0.hdslb.com/bfs/openplatform/202406/DVRww9ZG1719277949519.jpeg'
$arrExpo[16,0] = 17
$arrExpo[16,1] = '2024-07-20'
$arrExpo[16,2] = '南昌·漫拥动漫嘉年华Pro-追光启航'
$arrExpo[16,3] = '小蓝南路420号 洪州体育馆'
$arrExpo[16,4] = '2024.07.20 09:00-07.21 17:00'
$arrExpo[16,5] = 218
$arrExpo[16,6] = 52.5
$arrExpo[16,7] = 'https://show.bilibili.com/platform/detail.html?id=85796'
$arrExpo[16,8] = '//i1.hdslb.com/bfs/openplatform/202404/FawN3tPD1713364764414.png'
$arrExpo[17,0] = 18
$arrExpo[17,1] = '2024-07-21'
$arrExpo[17,2] = '乐平·CY境界次元动漫夏时庆'
$arrExpo[17,3] = '翥山西路182号 佳佳基大酒店'
$arrExpo[17,4] = '2024.07.21 10:00-07.21 17:00'
$arrExpo[17,5] = 142
$arrExpo[17,6] = 35
$arrExpo[17,7] = 'https://show.bilibili.com/platform/detail.html?id=86768'
$arrExpo[17,8] = '//i1.hdslb.com/bfs/openplatform/202406/3RWgXosx1717381178470.png'
$arrExpo[18,0] = 19
$arrExpo[18,1] = '2024-07-21'
$arrExpo[18,2] = '九江·SXD动漫嘉年华'
$arrExpo[18,3] = '湓浦街道大中路339号 百嘉洲际酒店'
$arrExpo[18,4] = '2024.07.21 10:00-07.21 17:30'
$arrExpo[18,5] = 63
$arrExpo[18,6] = 45
$arrExpo[18,7] = 'https://show.bilibili.com/platform/detail.html?id=86832'
$arrExpo[18,8] = '//i2.hdslb.com/bfs/openplatform/202406/Acs2Wqx71717394174913.jpeg'
$arrExpo[19,0] = 20
$arrExpo[19,1] = '2024-07-21'
$arrExpo[19,2] = '抚州·临次元08·盛夏动漫狂欢节'
$arrExpo[19,3] = '伍塘路1098号 乐课篮球公园'
$arrExpo[19,4] = '2024.07.21 10:00-07.21 16:00'
$arrExpo[19,5] = 80
$arrExpo[19,6] = 50
$arrExpo[19,7] = 'https://show.bilibili.com/platform/detail.html?id=87763'
$arrExpo[19,8] = '//i0.hdslb.com/bfs/openplatform/202406/6qgetbCh1718720523395.jpeg'
$arrExpo[20,0] = 21
$arrExpo[20,1] = '2024-07-21'
$arrExpo[20,2] = '萍乡·NL14动漫游戏展·夏日狂想曲'
$arrExpo[20,3] = '公园南路168号(近工行城北分理处) 梅生嘉华酒店'
$arrExpo[20,4] = '2024.07.21 10:00-07.21 17:00'
$arrExpo[20,5] = 78
$arrExpo[20,6] = 40
$arrExpo[20,7] = 'https://show.bilibili.com/platform/detail.html?id=86658'
$arrExpo[20,8] = '//i1.hdslb.com/bfs/openplatform/202405/bccpK1Zb1716969649865.jpeg'
$arrExpo[21,0] = 22
$arrExpo[21,1] = '2024-07-26'
$arrExpo[21,2] = '南昌·萌卡动漫展'
$arrExpo[21,3] = '八一桥街道青山南路118号蓝海购物广场F1 蓝海展览馆'
$arrExpo[21,4] = '2024.07.26 09:00-07.28 17:00'
$arrExpo[21,5] = 3070
$arrExpo[21,6] = 58.5
$arrExpo[21,7] = 'https://show.bilibili.com/platform/detail.html?id=86776'
$arrExpo[21,8] = '//i0.hdslb.com/bfs/openplatform/202406/WIQIJc741717410349369.jpeg'
$arrExpo[22,0] = 23
$arrExpo[22,1] = '2024-07-27'
$arrExpo[22,2] = '江西·次元星河动漫游戏嘉年华'
$arrExpo[22,3] = '九龙大道1177号 南昌绿地国际博览中心'
$arrExpo[22,4] = '2024.07.27 10:00-07.28 17:00'
$arrExpo[22,5] = 5456
$arrExpo[22,6] = 69
$arrExpo[22,7] = 'https://show.bilibili.com/platform/detail.html?id=85493'
$arrExpo[22,8] = '//i1.hdslb.com/bfs/openplatform/202405/jkKGgOqM1717141906659.png'
$arrExpo[23,0] = 24
$arrExpo[23,1] = '2024-07-27'
$arrExpo[23,2] = '赣州·马娘only'
$arrExpo[23,3] = '火车站广场正对面 赣州友尼宝国际酒店(赣州火车站店)'
$arrExpo[23,4] = '2024.07.27 09:00-07.27 17:00'
$arrExpo[23,5] = 33
$arrExpo[23,6] = 60
$arrExpo[23,7] = 'https://show.bilibili.com/platform/detail.html?id=86772'
$arrExpo[23,8] = '//i0.hdslb.com/bfs/openplatform/202406/BYe9CZzh1717172003064.png'
$arrExpo[24,0] = 25
$arrExpo[24,1] = '2024-07-28'
$arrExpo[24,2] = '赣州·明日方舟only叙拉古夜宴3.0暨同好交流茶话会'
$arrExpo[24,3] = '兴国路恒大帝景西门 江西长庚控股有限公司'
$arrExpo[24,4] = '2024.07.28 11:00-07.28 17:00'
$arrExpo[24,5] = 78
$arrExpo[24,6] = 56
$arrExpo[24,7] = 'https://show.bilibili.com/platform/detail.html?id=85688'
$arrExpo[24,8] = '//i1.hdslb.com/bfs/openplatform/202405/5AFwM8QV1715765287721.png'
$arrExpo[25,0] = 26
$arrExpo[25,1] = '2024-07-30'
$arrExpo[25,2] = '宜春·第三十五届静卿国风动漫文化展览会'
$arrExpo[25,3] = '宜阳大道19号(交通银行旁) 宜春安缦文华酒店'
$arrExpo[25,4] = '2024.07.30 09:00-07.30 17:00'
$arrExpo[25,5] = 500
$arrExpo[25,6] = 45
$arrExpo[25,7] = 'https://show.bilibili.com/platform/detail.html?id=86684'
$arrExpo[25,8] = '//i2.hdslb.com/bfs/openplatform/202406/meKBC0hU1719222126375.jpeg'
$arrExpo[26,0] = 27
$arrExpo[26,1] = '2024-07-31'
$arrExpo[26,2] = '万载·第八届馨缘动漫文化展'
$arrExpo[26,3] = '向阳路万载县幼儿园东南侧约60米 禧莱国际大酒店'
$arrExpo[26,4] = '2024.07.31 09:30-07.31 17:30'
$arrExpo[26,5] = 34
$arrExpo[26,6] = 40
$arrExpo[26,7] = 'https://show.bilibili.com/platform/detail.html?id=88294'
$arrExpo[26,8] = '//i2.hdslb.com/bfs/openplatform/202406/TQ4MvAGD1719195964738.jpeg'
$arrExpo[27,0] = 28
$arrExpo[27,1] = '2024-08-03'
$arrExpo[27,2] = '南昌·幻梦境国际动漫游戏嘉年华1th'
$arrExpo[27,3] = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$arrExpo[27,4] = '2024.08.03 09:00-08.04 17:30'
$arrExpo[27,5] = 3161
$arrExpo[27,6] = 64
$arrExpo[27,7] = 'https://show.bilibili.com/platform/detail.html?id=83980'
$arrExpo[27,8] = '//i0.hdslb.com/bfs/openplatform/202403/wRTbRtgD1710755902575.jpeg'
$arrExpo[28,0] = 29
$arrExpo[28,1] = '2024-08-03'
$arrExpo[28,2] = '吉安·COMIC LIFE周年庆典'
$arrExpo[28,3] = '东塘大道与阳明西路交叉路口往西约240米 吉安国际会展中心'
$arrExpo[28,4] = '2024.08.03 09:30-08.03 18:00'
$arrExpo[28,5] = 332
$arrExpo[28,6] = 52.1
$arrExpo[28,7] = 'https://show.bilibili.com/platform/detail.html?id=87164'
$arrExpo[28,8] = '//i1.hdslb.com/bfs/openplatform/202406/NWD9iQ9h1717598526259.jpeg'
$arrExpo[29,0] = 30
$arrExpo[29,1] = '2024-08-03'
$arrExpo[29,2] = '景德镇·第十五届瓷都ACG动漫游戏博览会'
$arrExpo[29,3] = '迎宾大道与寺山路交叉口东200米 陶博城'
$arrExpo[29,4] = '2024.08.03 09:00-08.04 17:00'
$arrExpo[29,5] = 2349
$arrExpo[29,6] = 55
$arrExpo[29,7] = 'https://show.bilibili.com/platform/detail.html?id=86341'
$arrExpo[29,8] = '//i0.hdslb.com/bfs/openplatform/202405/Wd6JiV3I1715953735690.png'
$arrExpo[30,0] = 31
$arrExpo[30,1] = '2024-08-03'
$arrExpo[30,2] = '景德镇·第十五届瓷都ACG动漫游戏博览会—马正阳内场票'
$arrExpo[30,3] = '迎宾大道与寺山路交叉口东200米 陶博城'
$arrExpo[30,4] = '2024.08.03 08:30-08.03 17:00'
$arrExpo[30,5] = 568
$arrExpo[30,6] = '已售罄'
$arrExpo[30,7] = 'https://show.bilibili.com/platform/detail.html?id=85981'
$arrExpo[30,8] = '//i2.hdslb.com/bfs/openplatform/202405/yevI9OGA1716445452947.png'
$arrExpo[31,0] = 32
$arrExpo[31,1] = '2024-08-03'
$arrExpo[31,2] = '樟树·第二届静卿国风动漫文化展览会'
$arrExpo[31,3] = '杏佛路89号 樟树银河国际酒店'
$arrExpo[31,4] = '2024.08.03 09:00-08.03 17:00'
$arrExpo[31,5] = 503
$arrExpo[31,6] = 45
$arrExpo[31,7] = 'https://show.bilibili.com/platform/detail.html?id=86683'
$arrExpo[31,8] = '//i2.hdslb.com/bfs/openplatform/202405/KD1hRj6P1716713054977.jpeg'
$arrExpo[32,0] = 33
$arrExpo[32,1] = '2024-08-03'
$arrExpo[32,2] = '萍乡·AU9夏至国漫展'
$arrExpo[32,3] = '安源中大道17号 壹号公馆（萍乡）'
$arrExpo[32,4] = '2024.08.03 10:00-08.03 17:00'
$arrExpo[32,5] = 102
$arrExpo[32,6] = 45
$arrExpo[32,7] = 'https://show.bilibili.com/platform/detail.html?id=86453'
$arrExpo[32,8] = '//i1.hdslb.com/bfs/openplatform/202406/hm1EACno1718936156944.jpeg'
$arrExpo[33,0] = 34
$arrExpo[33,1] = '2024-08-03'
$arrExpo[33,2] = '赣州·第一届环梦动漫游戏嘉年华'
$arrExpo[33,3] = '105国道东100米赣州毅德城国际会展中心 赣州毅德城国际会展中心'
$arrExpo[33,4] = '2024.08.03 09:00-08.05 17:00'
$arrExpo[33,5] = 161
$arrExpo[33,6] = 36.6
$arrExpo[33,7] = 'https://show.bilibili.com/platform/detail.html?id=87449'
$arrExpo[33,8] = '//i1.hdslb.com/bfs/openplatform/202406/iC3PuUfR1717740188790.jpeg'
$arrExpo[34,0] = 35
$arrExpo[34,1] = '2024-08-04'
$arrExpo[34,2] = '上饶·第十五届IX Group国风嘉年华暨十周年庆典'
$arrExpo[34,3] = '高铁经济试验区凤凰东大道16号7幢 上饶饶商金茂诚悦酒店(上饶高铁站)'
$arrExpo[34,4] = '2024.08.04 09:30-08.04 17:30'
$arrExpo[34,5] = 224
$arrExpo[34,6] = 60
$arrExpo[34,7] = 'https://show.bilibili.com/platform/detail.html?id=87225'
$arrExpo[34,8] = '//i2.hdslb.com/bfs/openplatform/202406/l5fIXZSX1717562269098.jpeg'
$arrExpo[35,0] = 36
$arrExpo[35,1] = '2024-08-04'
$arrExpo[35,2] = '九江·第一届异次元动漫嘉年华'
$arrExpo[35,3] = '长虹西大道兴城广场99号 九江半岛宾馆'
$arrExpo[35,4] = '2024.08.04 08:00-08.04 17:00'
$arrExpo[35,5] = 328
$arrExpo[35,6] = 45
$arrExpo[35,7] = 'https://show.bilibili.com/platform/detail.html?id=84407'
$arrExpo[35,8] = '//i2.hdslb.com/bfs/openplatform/202406/65hJjOfJ1717642614493.jpeg'
$arrExpo[36,0] = 37
$arrExpo[36,1] = '2024-08-05'
$arrExpo[36,2] = '上饶·囧喵喵国风动漫展'
$arrExpo[36,3] = '凤凰东大道与吴楚大道交叉路口北侧 饶派数字文创'
$arrExpo[36,4] = '2024.08.05 09:30-08.06 17:30'
$arrExpo[36,5] = 81
$arrExpo[36,6] = 65
$arrExpo[36,7] = 'https://show.bilibili.com/platform/detail.html?id=88050'
$arrExpo[36,8] = '//i1.hdslb.com/bfs/openplatform/202406/pTwNU5d41719229785696.jpeg'
$arrExpo[37,0] = 38
$arrExpo[37,1] = '2024-08-06'
$arrExpo[37,2] = '南昌·第一届异次元动漫嘉年华'
$arrExpo[37,3] = '民德路411号 东方豪景花园酒店(民德路店)'
$arrExpo[37,4] = '2024.08.06 08:00-08.06 17:00'
$arrExpo[37,5] = 485
$arrExpo[37,6] = 55
$arrExpo[37,7] = 'https://show.bilibili.com/platform/detail.html?id=84102'
$arrExpo[37,8] = '//i1.hdslb.com/bfs/openplatform/202405/BCA0owUW1716878997961.jpeg'
$arrExpo[38,0] = 39
$arrExpo[38,1] = '2024-08-08'
$arrExpo[38,2] = '赣州·第二届异次元动漫嘉年华'
$arrExpo[38,3] = '金辉路南3号大坪明德小学体育馆2层东侧201办公室 鲲伍体育·赣州经开区综合体育馆'
$arrExpo[38,4] = '2024.08.08 08:00-08.08 17:00'
$arrExpo[38,5] = 847
$arrExpo[38,6] = 45
$arrExpo[38,7] = 'https://show.bilibili.com/platform/detail.html?id=84184'
$arrExpo[38,8] = '//i1.hdslb.com/bfs/openplatform/202405/ayYIVKwP1716879335847.jpeg'
$arrExpo[39,0] = 40
$arrExpo[39,1] = '2024-08-10'
$arrExpo[39,2] = '南昌·花绒万兽第二聚'
$arrExpo[39,3] = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$arrExpo[39,4] = '2024.08.10 10:00-08.11 17:00'
$arrExpo[39,5] = 41
$arrExpo[39,6] = 188
$arrExpo[39,7] = 'https://show.bilibili.com/platform/detail.html?id=87600'
$arrExpo[39,8] = '//i1.hdslb.com/bfs/openplatform/202406/i0Ojsne01718693886054.png'
$arrExpo[40,0] = 41
$arrExpo[40,1] = '2024-08-10'
$arrExpo[40,2] = '吉安·WF无线次元新星动漫博览会'
$arrExpo[40,3] = '吉安南大道133号 吉安市全民健身中心'
$arrExpo[40,4] = '2024.08.10 09:00-08.10 17:00'
$arrExpo[40,5] = 23
$arrExpo[40,6] = 45
$arrExpo[40,7] = 'https://show.bilibili.com/platform/detail.html?id=88023'
$arrExpo[40,8] = '//i0.hdslb.com/bfs/openplatform/202406/f95zVAmw1718246635629.jpeg'
$arrExpo[41,0] = 42
$arrExpo[41,1] = '2024-08-10'
$arrExpo[41,2] = '高安·第二届静卿国风动漫文化展览会'
$arrExpo[41,3] = '华林中路606号 高安华鼎国际大酒店'
$arrExpo[41,4] = '2024.08.10 09:00-08.10 17:00'
$arrExpo[41,5] = 458
$arrExpo[41,6] = 45
$arrExpo[41,7] = 'https://show.bilibili.com/platform/detail.html?id=86682'
$arrExpo[41,8] = '//i2.hdslb.com/bfs/openplatform/202405/UwvNYGne1716711642772.jpeg'
$arrExpo[42,0] = 43
$arrExpo[42,1] = '2024-08-15'
$arrExpo[42,2] = '上饶·次元重现夏日嘉年华'
$arrExpo[42,3] = '普济巷地委大院北侧约90米 四季体育运动馆'
$arrExpo[42,4] = '2024.08.15 09:30-08.15 17:30'
$arrExpo[42,5] = 52
$arrExpo[42,6] = 58
$arrExpo[42,7] = 'https://show.bilibili.com/platform/detail.html?id=87679'
$arrExpo[42,8] = '//i0.hdslb.com/bfs/openplatform/202406/fxlKV2SL1718784421064.jpeg'
$arrExpo[43,0] = 44
$arrExpo[43,1] = '2024-08-24'
$arrExpo[43,2] = '南昌·第四届龙年动漫展——暑假最后的狂欢'
$arrExpo[43,3] = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$arrExpo[43,4] = '2024.08.24 10:00-08.25 18:00'
$arrExpo[43,5] = 521
$arrExpo[43,6] = 45
$arrExpo[43,7] = 'https://show.bilibili.com/platform/detail.html?id=87135'
$arrExpo[43,8] = '//i0.hdslb.com/bfs/openplatform/202406/mDtqZeQd1718033555304.jpeg'
$wsExpo.Range("A2:I45").Value = $arrExpo

$wsExpo.Range("A1").Copy()
$wsExpo.Range("A2:A45").PasteSpecial(-4122)
$wsExpo.Range("B2:B45").Style = "Normal"

# ---- Sheet: 演出 (Performance) -- append data row 2 ----
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("B2").NumberFormat = "@"

$arrShow = New-Object 'object[,]' 1,9
$arrShow[0,0] = 1
$arrShow[0,1] = '2024-07-13'
$arrShow[0,2] = '江西·东方LiveParty×THO03幻想Strawberry~！！'
$arrShow[0,3] = '上海路543号520Park文创公园21号01区域 瓦肆VAS NANCHANG'
$arrShow[0,4] = '2024.07.13 20:30-07.13 23:00'
$arrShow[0,5] = 84
$arrShow[0,6] = 160
$arrShow[0,7] = 'https://show.bilibili.com/platform/detail.html?id=87366'
$arrShow[0,8] = '//i0.hdslb.com/bfs/openplatform/202406/1L3I6Qmg1718292516616.jpeg'
$wsShow.Range("A2:I2").Value = $arrShow

$wsShow.Range("A1").Copy()
$wsShow.Range("A2").PasteSpecial(-4122)
$wsShow.Range("B2").Style = "Normal"

# ---- Sheet: 全部类型 (All types) -- correct F column (想去人数), rows 2..46 ----
$wsAll = $wb.Worksheets.Item("全部类型")
$arrAllF = New-Object 'object[,]' 45,1
$arrAllF[0,0] = 1862
$arrAllF[1,0] = 277
$arrAllF[2,0] = 32
$arrAllF[3,0] = 62
$arrAllF[4,0] = 44
$arrAllF[5,0] = 169
$arrAllF[6,0] = 3700
$arrAllF[7,0] = 165
$arrAllF[8,0] = 108
$arrAllF[9,0] = 84
$arrAllF[10,0] = 84
$arrAllF[11,0] = 67
$arrAllF[12,0] = 72
$arrAllF[13,0] = 646
$arrAllF[14,0] = 146
$arrAllF[15,0] = 841
$arrAllF[16,0] = 41
$arrAllF[17,0] = 218
$arrAllF[18,0] = 142
$arrAllF[19,0] = 63
$arrAllF[20,0] = 80
$arrAllF[21,0] = 78
$arrAllF[22,0] = 3070
$arrAllF[23,0] = 5456
$arrAllF[24,0] = 33
$arrAllF[25,0] = 78
$arrAllF[26,0] = 500
$arrAllF[27,0] = 34
$arrAllF[28,0] = 3161
$arrAllF[29,0] = 332
$arrAllF[30,0] = 2349
$arrAllF[31,0] = 568
$arrAllF[32,0] = 503
$arrAllF[33,0] = 102
$arrAllF[34,0] = 161
$arrAllF[35,0] = 224
$arrAllF[36,0] = 328
$arrAllF[37,0] = 81
$arrAllF[38,0] = 485
$arrAllF[39,0] = 847
$arrAllF[40,0] = 41
$arrAllF[41,0] = 23
$arrAllF[42,0] = 458
$arrAllF[43,0] = 52
$arrAllF[44,0] = 521
$wsAll.Range("F2:F46").Value = $arrAllF

$wsExpo.Select() | Out-Null
